$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header renames ---
$ws.Range("A1").Value = "industry"
$ws.Range("B1").Value = "unit"
$ws.Range("C1").Value = "process"
$ws.Range("D1").Value = "carbon (kg CO2 eq)"
$ws.Range("E1").Value = "ced (MJ)"
$ws.Range("F1").Value = "climate change (kg CO2 eq)"
$ws.Range("G1").Value = "region"

# --- Data updates (D/E/F columns, rows 2-166) ---
$ws.Range("D2").Value = 1.2159612
$ws.Range("E2").Value = 15.486373
$ws.Range("F2").Value = 0.000033904224
$ws.Range("D3").Value = 4.813070866666667
$ws.Range("E3").Value = 71.689488
$ws.Range("F3").Value = 0.00013420119
$ws.Range("D4").Value = 0.7504491333333334
$ws.Range("E4").Value = 11.287101
$ws.Range("F4").Value = 0.000020924513
$ws.Range("D5").Value = 3.959920333333333
$ws.Range("E5").Value = 59.004986
$ws.Range("F5").Value = 0.00011041309
$ws.Range("D6").Value = 11.99553666666667
$ws.Range("E6").Value = 142.93184
$ws.Range("F6").Value = 0.0003344674
$ws.Range("D7").Value = 3.818
$ws.Range("E7").Value = 78.28485
$ws.Range("F7").Value = 0.00010645597
$ws.Range("D8").Value = 21.70634933333334
$ws.Range("E8").Value = 36.68254
$ws.Range("F8").Value = 0.00060523061
$ws.Range("D9").Value = 37.20312466666667
$ws.Range("E9").Value = 194.53125
$ws.Range("F9").Value = 0.0010373218
$ws.Range("D10").Value = 3.4506146
$ws.Range("E10").Value = 57.44649
$ws.Range("F10").Value = 0.00009621229
$ws.Range("D11").Value = 0.46758162
$ws.Range("E11").Value = 8.3697714
$ws.Range("F11").Value = 0.000013037416
$ws.Range("D12").Value = 2.496044
$ws.Range("E12").Value = 41.74194
$ws.Range("F12").Value = 0.000069596331
$ws.Range("D13").Value = 170.015
$ws.Range("E13").Value = 139.48
$ws.Range("F13").Value = 0.0047404693
$ws.Range("D14").Value = 47790.3
$ws.Range("E14").Value = 92660.0
$ws.Range("F14").Value = 1.3325204
$ws.Range("D15").Value = 750.6249333333334
$ws.Range("E15").Value = 9903.7414
$ws.Range("F15").Value = 0.020929415
$ws.Range("D16").Value = 36500.778
$ws.Range("E16").Value = 72798.498
$ws.Range("F16").Value = 1.0177385
$ws.Range("D17").Value = 117.52
$ws.Range("E17").Value = 389.7
$ws.Range("F17").Value = 0.0032767694
$ws.Range("D18").Value = 1.892
$ws.Range("E18").Value = 30.966238
$ws.Range("F18").Value = 0.00005275398
$ws.Range("D19").Value = 0.6074231533333334
$ws.Range("E19").Value = 9.0839967
$ws.Range("F19").Value = 0.000016936569
$ws.Range("D20").Value = 0.8643385333333333
$ws.Range("E20").Value = 13.460445
$ws.Range("F20").Value = 0.000024100051
$ws.Range("D21").Value = 79.29
$ws.Range("E21").Value = 95.53
$ws.Range("F21").Value = 0.0022108156
$ws.Range("D22").Value = 26.47878733333334
$ws.Range("E22").Value = 392.3316
$ws.Range("F22").Value = 0.00073829885
$ws.Range("D23").Value = 0.4933606866666667
$ws.Range("E23").Value = 7.3939529
$ws.Range("F23").Value = 0.000013756205
$ws.Range("D24").Value = 24.65980733333333
$ws.Range("E24").Value = 365.38597
$ws.Range("F24").Value = 0.00068758087
$ws.Range("D25").Value = 4.9818182
$ws.Range("E25").Value = 7.1424762
$ws.Range("F25").Value = 0.00013890631
$ws.Range("D26").Value = 11.55149
$ws.Range("E26").Value = 0.972
$ws.Range("F26").Value = 0.00032208619
$ws.Range("D27").Value = 15.51872733333333
$ws.Range("E27").Value = 42.572344
$ws.Range("F27").Value = 0.00043270329
$ws.Range("D28").Value = 13.095508
$ws.Range("E28").Value = 251.41399
$ws.Range("F28").Value = 0.00036513752
$ws.Range("D29").Value = 0.5422446000000001
$ws.Range("E29").Value = 8.1182574
$ws.Range("F29").Value = 0.000015119218
$ws.Range("D30").Value = 11.086986
$ws.Range("E30").Value = 212.48668
$ws.Range("F30").Value = 0.00030913459
$ws.Range("D31").Value = 23700.0
$ws.Range("E31").Value = 363300.0
$ws.Range("F31").Value = 0.66081889
$ws.Range("D32").Value = 311.52784
$ws.Range("E32").Value = 4110.2966
$ws.Range("F32").Value = 0.0086862228
$ws.Range("D33").Value = 15981.804
$ws.Range("E33").Value = 244767.4
$ws.Range("F33").Value = 0.44561511
$ws.Range("D34").Value = 33300.0
$ws.Range("E34").Value = 465150.0
$ws.Range("F34").Value = 0.92849236
$ws.Range("D35").Value = 1100.0538
$ws.Range("E35").Value = 14514.104
$ws.Range("F35").Value = 0.030672419
$ws.Range("D36").Value = 25250.01333333333
$ws.Range("E36").Value = 352491.03
$ws.Range("F36").Value = 0.70403737
$ws.Range("D37").Value = 34900.0
$ws.Range("E37").Value = 481950.0
$ws.Range("F37").Value = 0.9731046
$ws.Range("D38").Value = 3106.034266666667
$ws.Range("E38").Value = 40980.999
$ws.Range("F38").Value = 0.086604477
$ws.Range("D39").Value = 23454.17266666667
$ws.Range("E39").Value = 323201.16
$ws.Range("F39").Value = 0.65396456
$ws.Range("D40").Value = 5.458457933333333
$ws.Range("E40").Value = 57.658152
$ws.Range("F40").Value = 0.00015219629
$ws.Range("D41").Value = 449.3869
$ws.Range("E41").Value = 801.26
$ws.Range("F41").Value = 0.012530099
$ws.Range("D42").Value = 14.373962
$ws.Range("E42").Value = 171.27177
$ws.Range("F42").Value = 0.00040078421
$ws.Range("D43").Value = 366.73444
$ws.Range("E43").Value = 681.56224
$ws.Range("F43").Value = 0.010225529
$ws.Range("D44").Value = 144.5000533333333
$ws.Range("E44").Value = 109.92602
$ws.Range("F44").Value = 0.0040290449
$ws.Range("D45").Value = 44.537884
$ws.Range("E45").Value = 530.68752
$ws.Range("F45").Value = 0.0012418344
$ws.Range("D46").Value = 10.14680666666667
$ws.Range("E46").Value = 15.32245
$ws.Range("F46").Value = 0.00028291988
$ws.Range("D47").Value = 28.804014
$ws.Range("E47").Value = 442.65967
$ws.Range("F47").Value = 0.00080313234
$ws.Range("D48").Value = 19.02690333333334
$ws.Range("E48").Value = 343.49938
$ws.Range("F48").Value = 0.00053052056
$ws.Range("D49").Value = 28.70624266666667
$ws.Range("E49").Value = 441.66806
$ws.Range("F49").Value = 0.00080040622
$ws.Range("D50").Value = 1.683935266666667
$ws.Range("E50").Value = 32.958787
$ws.Range("F50").Value = 0.000046952584
$ws.Range("D51").Value = 17.105734
$ws.Range("E51").Value = 314.69843
$ws.Range("F51").Value = 0.00047695324
$ws.Range("D52").Value = 16.93467666666666
$ws.Range("E52").Value = 311.55145
$ws.Range("F52").Value = 0.00047218371
$ws.Range("D53").Value = 3.818
$ws.Range("E53").Value = 78.28485
$ws.Range("F53").Value = 0.00010645597
$ws.Range("D54").Value = 0.9325907333333333
$ws.Range("E54").Value = 13.918732
$ws.Range("F54").Value = 0.000026003105
$ws.Range("D55").Value = 3.414042733333333
$ws.Range("E55").Value = 69.273593
$ws.Range("F55").Value = 0.00009519257
$ws.Range("D56").Value = 2.87
$ws.Range("E56").Value = 56.21302
$ws.Range("F56").Value = 0.000080023215
$ws.Range("D57").Value = 0.036623172
$ws.Range("E57").Value = 0.66124081
$ws.Range("F57").Value = 0.0000010211512
$ws.Range("D58").Value = 3.8906146
$ws.Range("E58").Value = 57.44649
$ws.Range("F58").Value = 0.00010848066
$ws.Range("D59").Value = 3.132987933333333
$ws.Range("E59").Value = 32.252174
$ws.Range("F59").Value = 0.000087356017
$ws.Range("D60").Value = 558.68876
$ws.Range("E60").Value = 8378.7135
$ws.Range("F60").Value = 0.015577725
$ws.Range("D61").Value = 1.892
$ws.Range("E61").Value = 30.966238
$ws.Range("F61").Value = 0.00005275398
$ws.Range("D62").Value = 50.03605333333334
$ws.Range("E62").Value = 324.76761
$ws.Range("F62").Value = 0.0013951379
$ws.Range("D63").Value = 310.2294
$ws.Range("E63").Value = 3696.5131
$ws.Range("F63").Value = 0.0086500188
$ws.Range("D64").Value = 0.18
$ws.Range("E64").Value = 3.85378
$ws.Range("F64").Value = 0.0000050188776
$ws.Range("D65").Value = 798.9441333333334
$ws.Range("E65").Value = 9519.7533
$ws.Range("F65").Value = 0.022276682
$ws.Range("D66").Value = 4.817188666666667
$ws.Range("E66").Value = 51.381532
$ws.Range("F66").Value = 0.000134316
$ws.Range("D67").Value = 357.8199266666667
$ws.Range("E67").Value = 824.52574
$ws.Range("F67").Value = 0.009976969
$ws.Range("D68").Value = 111.2689466666667
$ws.Range("E68").Value = 1325.816
$ws.Range("F68").Value = 0.0031024734
$ws.Range("D69").Value = 25.95697333333333
$ws.Range("E69").Value = 309.70939
$ws.Range("F69").Value = 0.00072374928
$ws.Range("D70").Value = 121.14674
$ws.Range("E70").Value = 1443.9395
$ws.Range("F70").Value = 0.0033778925
$ws.Range("D71").Value = 98.81022
$ws.Range("E71").Value = 1177.7905
$ws.Range("F71").Value = 0.0027550912
$ws.Range("D72").Value = 801.5832
$ws.Range("E72").Value = 9551.6248
$ws.Range("F72").Value = 0.022350267
$ws.Range("D73").Value = 94.67383333333333
$ws.Range("E73").Value = 1128.5037
$ws.Range("F73").Value = 0.0026397576
$ws.Range("D74").Value = 22.39038
$ws.Range("E74").Value = 267.21616
$ws.Range("F74").Value = 0.00062430321
$ws.Range("D75").Value = 26.580362
$ws.Range("E75").Value = 317.13872
$ws.Range("F75").Value = 0.00074113102
$ws.Range("D76").Value = 35.73024466666666
$ws.Range("E76").Value = 426.16622
$ws.Range("F76").Value = 0.00099625402
$ws.Range("D77").Value = 39.039358
$ws.Range("E77").Value = 465.59569
$ws.Range("F77").Value = 0.0010885209
$ws.Range("D78").Value = 120.11264
$ws.Range("E78").Value = 1431.6178
$ws.Range("F78").Value = 0.0033490591
$ws.Range("D79").Value = 8.870582666666667
$ws.Range("E79").Value = 106.12212
$ws.Range("F79").Value = 0.00024733539
$ws.Range("D80").Value = 602.0023066666668
$ws.Range("E80").Value = 7173.5347
$ws.Range("F80").Value = 0.016785422
$ws.Range("D81").Value = 253.51128
$ws.Range("E81").Value = 3021.1184
$ws.Range("F81").Value = 0.0070685672
$ws.Range("D82").Value = 30.663164
$ws.Range("E82").Value = 365.78984
$ws.Range("F82").Value = 0.00085497038
$ws.Range("D83").Value = 7.796
$ws.Range("E83").Value = 79.677098
$ws.Range("F83").Value = 0.00021737316
$ws.Range("D84").Value = 7.967754
$ws.Range("E84").Value = 83.925865
$ws.Range("F84").Value = 0.00022216212
$ws.Range("D85").Value = 8.104398
$ws.Range("E85").Value = 86.385279
$ws.Range("F85").Value = 0.00022597213
$ws.Range("D86").Value = 8.081428666666667
$ws.Range("E86").Value = 86.086044
$ws.Range("F86").Value = 0.00022533168
$ws.Range("D87").Value = 7.943595333333334
$ws.Range("E87").Value = 82.889989
$ws.Range("F87").Value = 0.00022148851
$ws.Range("D88").Value = 8.270162
$ws.Range("E88").Value = 87.06109
$ws.Range("F88").Value = 0.00023059407
$ws.Range("D89").Value = 8.404252666666668
$ws.Range("E89").Value = 90.424835
$ws.Range("F89").Value = 0.00023433286
$ws.Range("D90").Value = 8.849138666666667
$ws.Range("E90").Value = 96.770762
$ws.Range("F90").Value = 0.00024673747
$ws.Range("D91").Value = 8.163238666666668
$ws.Range("E91").Value = 85.407134
$ws.Range("F91").Value = 0.00022761275
$ws.Range("D92").Value = 8.186344
$ws.Range("E92").Value = 85.764122
$ws.Range("F92").Value = 0.000228257
$ws.Range("D93").Value = 8.073382666666667
$ws.Range("E93").Value = 83.478655
$ws.Range("F93").Value = 0.00022510732
$ws.Range("D94").Value = 8.236392
$ws.Range("E94").Value = 86.190674
$ws.Range("F94").Value = 0.00022965247
$ws.Range("D95").Value = 8.165737333333334
$ws.Range("E95").Value = 85.385908
$ws.Range("F95").Value = 0.00022768243
$ws.Range("D96").Value = 8.189875333333335
$ws.Range("E96").Value = 88.972062
$ws.Range("F96").Value = 0.00022835546
$ws.Range("D97").Value = 7.875716000000001
$ws.Range("E97").Value = 82.866194
$ws.Range("F97").Value = 0.00021959585
$ws.Range("D98").Value = 8.562822
$ws.Range("E98").Value = 92.26768
$ws.Range("F98").Value = 0.0002387542
$ws.Range("D99").Value = 8.737582666666668
$ws.Range("E99").Value = 96.807237
$ws.Range("F99").Value = 0.000243627
$ws.Range("D100").Value = 7.766506
$ws.Range("E100").Value = 81.080848
$ws.Range("F100").Value = 0.0002165508
$ws.Range("D101").Value = 7.767872666666666
$ws.Range("E101").Value = 81.376162
$ws.Range("F101").Value = 0.0002165889
$ws.Range("D102").Value = 7.968067333333333
$ws.Range("E102").Value = 83.306462
$ws.Range("F102").Value = 0.00022217087
$ws.Range("D103").Value = 7.683113333333333
$ws.Range("E103").Value = 81.015131
$ws.Range("F103").Value = 0.00021422559
$ws.Range("D104").Value = 3.5209
$ws.Range("E104").Value = 57.261
$ws.Range("F104").Value = 0.000098172034
$ws.Range("D105").Value = 2.6838
$ws.Range("E105").Value = 46.688
$ws.Range("F105").Value = 0.000074831465
$ws.Range("D106").Value = 3.8207
$ws.Range("E106").Value = 70.414
$ws.Range("F106").Value = 0.00010653125
$ws.Range("D107").Value = 2.496044
$ws.Range("E107").Value = 41.74194
$ws.Range("F107").Value = 0.000069596331
$ws.Range("D108").Value = 3.0338722
$ws.Range("E108").Value = 44.958568
$ws.Range("F108").Value = 0.000084592406
$ws.Range("D109").Value = 3.039791666666666
$ws.Range("E109").Value = 47.556399
$ws.Range("F109").Value = 0.000084757458
$ws.Range("D110").Value = 3.541340933333333
$ws.Range("E110").Value = 60.877152
$ws.Range("F110").Value = 0.000098741981
$ws.Range("D111").Value = 4.397107
$ws.Range("E111").Value = 80.517389
$ws.Range("F111").Value = 0.00012260301
$ws.Range("D112").Value = 6.472074733333334
$ws.Range("E112").Value = 118.96978
$ws.Range("F112").Value = 0.00018045862
$ws.Range("D113").Value = 3.177242733333334
$ws.Range("E113").Value = 42.465158
$ws.Range("F113").Value = 0.000088589959
$ws.Range("D114").Value = 3.279198733333333
$ws.Range("E114").Value = 42.163446
$ws.Range("F114").Value = 0.000091432761
$ws.Range("D115").Value = 2.804837533333334
$ws.Range("E115").Value = 48.406753
$ws.Range("F115").Value = 0.000078206313
$ws.Range("D116").Value = 2.942537333333334
$ws.Range("E116").Value = 52.536501
$ws.Range("F116").Value = 0.000082045749
$ws.Range("D117").Value = 3.006797266666667
$ws.Range("E117").Value = 54.463717
$ws.Range("F117").Value = 0.000083837485
$ws.Range("D118").Value = 3.0343372
$ws.Range("E118").Value = 55.289667
$ws.Range("F118").Value = 0.000084605372
$ws.Range("D119").Value = 3.009789866666667
$ws.Range("E119").Value = 54.773674
$ws.Range("F119").Value = 0.000083920928
$ws.Range("D120").Value = 3.197522
$ws.Range("E120").Value = 47.537062
$ws.Range("F120").Value = 0.000089155397
$ws.Range("D121").Value = 3.670218133333333
$ws.Range("E121").Value = 56.661292
$ws.Range("F121").Value = 0.00010233542
$ws.Range("D122").Value = 3.563410333333333
$ws.Range("E122").Value = 60.70187
$ws.Range("F122").Value = 0.000099357335
$ws.Range("D123").Value = 3.432213933333333
$ws.Range("E123").Value = 41.635056
$ws.Range("F123").Value = 0.000095699232
$ws.Range("D124").Value = 3.5852292
$ws.Range("E124").Value = 41.106667
$ws.Range("F124").Value = 0.000099965703
$ws.Range("D125").Value = 3.0653736
$ws.Range("E125").Value = 43.793828
$ws.Range("F125").Value = 0.000085470749
$ws.Range("D126").Value = 2.804837533333334
$ws.Range("E126").Value = 48.406753
$ws.Range("F126").Value = 0.000078206313
$ws.Range("D127").Value = 3.024236333333334
$ws.Range("E127").Value = 54.491758
$ws.Range("F127").Value = 0.000084323734
$ws.Range("D128").Value = 3.0343372
$ws.Range("E128").Value = 55.289667
$ws.Range("F128").Value = 0.000084605372
$ws.Range("D129").Value = 23.40999866666667
$ws.Range("E129").Value = 343.21701
$ws.Range("F129").Value = 0.00065273288
$ws.Range("D130").Value = 24.87992133333334
$ws.Range("E130").Value = 367.84722
$ws.Range("F130").Value = 0.00069371823
$ws.Range("D131").Value = 23.447178
$ws.Range("E131").Value = 345.762
$ws.Range("F131").Value = 0.00065376954
$ws.Range("D132").Value = 23.660106
$ws.Range("E132").Value = 347.89377
$ws.Range("F132").Value = 0.00065970654
$ws.Range("D133").Value = 23.15372133333333
$ws.Range("E133").Value = 340.15185
$ws.Range("F133").Value = 0.00064558719
$ws.Range("D134").Value = 23.42401266666667
$ws.Range("E134").Value = 343.89384
$ws.Range("F134").Value = 0.00065312363
$ws.Range("D135").Value = 24.35832933333333
$ws.Range("E135").Value = 359.89726
$ws.Range("F135").Value = 0.00067917485
$ws.Range("D136").Value = 23.94754
$ws.Range("E136").Value = 355.85996
$ws.Range("F136").Value = 0.00066772097
$ws.Range("D137").Value = 24.82088733333334
$ws.Range("E137").Value = 366.77249
$ws.Range("F137").Value = 0.0006920722
$ws.Range("D138").Value = 23.77793533333333
$ws.Range("E138").Value = 350.07658
$ws.Range("F138").Value = 0.00066299193
$ws.Range("D139").Value = 11.49915533333333
$ws.Range("E139").Value = 214.1629
$ws.Range("F139").Value = 0.00032062696
$ws.Range("D140").Value = 9.485888666666668
$ws.Range("E140").Value = 134.97303
$ws.Range("F140").Value = 0.00026449174
$ws.Range("D141").Value = 10.07873266666667
$ws.Range("E141").Value = 180.33996
$ws.Range("F141").Value = 0.00028102182
$ws.Range("D142").Value = 11.086986
$ws.Range("E142").Value = 212.48668
$ws.Range("F142").Value = 0.00030913459
$ws.Range("D143").Value = 8.282163333333335
$ws.Range("E143").Value = 121.78768
$ws.Range("F143").Value = 0.00023092868
$ws.Range("D144").Value = 13.71599266666667
$ws.Range("E144").Value = 184.81033
$ws.Range("F144").Value = 0.00038243826
$ws.Range("D145").Value = 18.26592466666667
$ws.Range("E145").Value = 183.45407
$ws.Range("F145").Value = 0.00050930244
$ws.Range("D146").Value = 14.088858
$ws.Range("E146").Value = 188.28705
$ws.Range("F146").Value = 0.00039283474
$ws.Range("D147").Value = 10.57226933333333
$ws.Range("E147").Value = 189.99296
$ws.Range("F147").Value = 0.00029478291
$ws.Range("D148").Value = 8.714413333333333
$ws.Range("E148").Value = 162.41603
$ws.Range("F148").Value = 0.00024298097
$ws.Range("D149").Value = 7.443728
$ws.Range("E149").Value = 126.51838
$ws.Range("F149").Value = 0.00020755089
$ws.Range("D150").Value = 13.483894
$ws.Range("E150").Value = 167.47214
$ws.Range("F150").Value = 0.00037596674
$ws.Range("D151").Value = 10.09811
$ws.Range("E151").Value = 173.77009
$ws.Range("F151").Value = 0.00028156209
$ws.Range("D152").Value = 36.04956533333333
$ws.Range("E152").Value = 526.6326
$ws.Range("F152").Value = 0.0010051575
$ws.Range("D153").Value = 35.83358666666667
$ws.Range("E153").Value = 528.46667
$ws.Range("F153").Value = 0.00099913548
$ws.Range("D154").Value = 34.399242
$ws.Range("E154").Value = 500.07364
$ws.Range("F154").Value = 0.00095914213
$ws.Range("D155").Value = 4.495870133333334
$ws.Range("E155").Value = 83.771835
$ws.Range("F155").Value = 0.00012535679
$ws.Range("D156").Value = 4.440957133333334
$ws.Range("E156").Value = 82.853513
$ws.Range("F156").Value = 0.00012382567
$ws.Range("D157").Value = 4.459677466666666
$ws.Range("E157").Value = 83.166577
$ws.Range("F157").Value = 0.00012434764
$ws.Range("D158").Value = 4.501977200000001
$ws.Range("E158").Value = 85.364631
$ws.Range("F158").Value = 0.00012552707
$ws.Range("D159").Value = 5.0251234
$ws.Range("E159").Value = 84.14837
$ws.Range("F159").Value = 0.00014011377
$ws.Range("D160").Value = 7.111703333333334
$ws.Range("E160").Value = 114.06269
$ws.Range("F160").Value = 0.00019829315
$ws.Range("D161").Value = 18.30489266666667
$ws.Range("E161").Value = 237.68303
$ws.Range("F161").Value = 0.00051038897
$ws.Range("D162").Value = 4.686698866666667
$ws.Range("E162").Value = 79.97873
$ws.Range("F162").Value = 0.0001306776
$ws.Range("D163").Value = 21.43328066666667
$ws.Range("E163").Value = 271.371
$ws.Range("F163").Value = 0.00059761673
$ws.Range("D164").Value = 4.844162133333334
$ws.Range("E164").Value = 81.676307
$ws.Range("F164").Value = 0.00013506809
$ws.Range("D165").Value = 72.47620666666667
$ws.Range("E165").Value = 737.844
$ws.Range("F165").Value = 0.0020208289
$ws.Range("D166").Value = 7.673719333333334
$ws.Range("E166").Value = 119.42972
$ws.Range("F166").Value = 0.00021396365

# --- Header cell comments ---
$ws.Range("A1").AddComment("Data type: Categorical (text)")
$ws.Range("B1").AddComment("Data type: Various (e.g. kg, kWh)")
$ws.Range("C1").AddComment("Data type: Categorical (text)")
$ws.Range("D1").AddComment("Data type: Carbon footprint")
$ws.Range("E1").AddComment("Data type: Cumulative energy demand")
$ws.Range("F1").AddComment("Data type: Climate change impact")
$ws.Range("G1").AddComment("Data type: Categorical (text)")
